$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the D62:D67 formulas as one range write so the engine collapses
# them into a shared-formula group (matching si="1" in the target diff).
$ws.Range("D62:D67").Formula = "=C62/(24*60)"

# --- New row 68: 四方坪站 (2025-10-04) ---
$ws.Range("A68").Value = 45934
$ws.Range("B68").Value = "四方坪站"
$ws.Range("C68").Formula = "=14070/127"
$ws.Range("D68").Formula = "=C68/(24*60)"
$ws.Range("E68").Formula = "=8238.48/127"
$ws.Range("F68").Formula = "=2810.88/127"
$ws.Range("G68").Formula = "=8238.48/(14070/60)"
$ws.Range("H68").Formula = "=354/127"

# --- New row 69: 高岭站 (2025-10-04) ---
$ws.Range("A69").Value = 45934
$ws.Range("B69").Value = "高岭站"
$ws.Range("C69").Formula = "=5345/36"
$ws.Range("D69").Formula = "=C69/(24*60)"
$ws.Range("E69").Formula = "=3987.85/36"
$ws.Range("F69").Formula = "=964.79/36"
$ws.Range("G69").Formula = "=3987.85/(5345/60)"
$ws.Range("H69").Formula = "=141/36"

$ws.Range("J70").Select()
